$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("G2").Value = 8.19
$ws.Range("G3").Value = 16.38
$ws.Range("G4").Value = 16.38
$ws.Range("G5").Value = 32.76
$ws.Range("G6").Value = 16.38

$ws.Range("F10").Value = 0.24305555555555555
$ws.Range("G10").Value = 26.39
$ws.Range("F11").Value = 0.1875

$ws.Activate()
$ws.Range("M17").Select()
